# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values computed for rows 2-15
$kValues = @{
    2  = 1
    3  = 8
    4  = 4
    5  = 5
    6  = 2
    7  = 2
    8  = 3
    9  = 3
    10 = 3
    11 = 0
    12 = 3
    13 = 6
    14 = 1
    15 = 7
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
